# Daily attendance processing - 2026-02-01 16:45:09
# Normalize the "Recorded By" (column G) values: several rows list the
# recording users in an inconsistent order (e.g. "System" should be listed
# after the other recorder rather than before it). This pass walks the
# "Recorded By" column and re-orders a small set of known value patterns
# to their corrected form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Known "Recorded By" strings -> their corrected (re-ordered) form.
$map = @{
    "System, system, backup@backdoor.com" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
